$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @("D2", "27.370.98"),
    @("E2", "  -0.64%  "),
    @("D3", "1.790.85"),
    @("E3", "  -1.73%  "),
    @("D4", "1.002"),
    @("E4", "  -0.34%  "),
    @("D5", "340.57"),
    @("E5", "  +1.02%  "),
    @("D6", "0.9993"),
    @("E6", "  -0.14%  "),
    @("D7", "0.3946"),
    @("E7", "  +3.16%  "),
    @("D8", "0.3468"),
    @("E8", "  -1.70%  "),
    @("D9", "48.16"),
    @("E9", "  -3.50%  "),
    @("D10", "1.199"),
    @("E10", "  -2.87%  "),
    @("D11", "0.07509"),
    @("E11", "  -2.71%  "),
    @("D12", "0.9996"),
    @("E12", "  -0.20%  "),
    @("D13", "21.84"),
    @("E13", "  -1.83%  "),
    @("D14", "6.503"),
    @("E14", "  -1.70%  "),
    @("D15", "1.786.36"),
    @("E15", "  -2.24%  "),
    @("D16", "7.135"),
    @("E16", "  -0.60%  "),
    @("D17", "0.00001098"),
    @("E17", "  -2.16%  "),
    @("D18", "0.06695"),
    @("E18", "  -0.35%  "),
    @("D19", "84.97"),
    @("E19", "  -2.49%  "),
    @("D20", "0.9990"),
    @("E20", "  -0.12%  "),
    @("D21", "17.78"),
    @("E21", "  +0.89%  "),
    @("D22", "6.528"),
    @("E22", "  -0.09%  "),
    @("D23", "27.364.22"),
    @("E23", "  -0.77%  "),
    @("D24", "12.44"),
    @("E24", "  -5.24%  "),
    @("D25", "2.411"),
    @("E25", "  -2.70%  "),
    @("E26", "  -3.72%  "),
    @("D27", "2.508"),
    @("E27", "  -5.59%  "),
    @("D28", "1.461"),
    @("E28", "  -1.22%  "),
    @("D29", "157.98"),
    @("E29", "  +3.54%  "),
    @("D30", "1.988.99"),
    @("E30", "  -2.08%  "),
    @("D31", "136.49"),
    @("E31", "  +1.07%  "),
    @("D32", "4.033"),
    @("E32", "  -1.30%  "),
    @("D33", "6.021"),
    @("E33", "  -4.53%  "),
    @("E34", "  +0.27%  "),
    @("D35", "13.05"),
    @("E35", "  -5.78%  "),
    @("B36", "WEMIXTOKEN"),
    @("C36", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"),
    @("D36", "1.623"),
    @("E36", "  -4.46%  "),
    @("B37", "Hedera"),
    @("C37", "https://coinranking.com/coin/jad286TjB+hedera-hbar"),
    @("D37", "0.06559"),
    @("E37", "  +0.84%  "),
    @("D38", "0.02428"),
    @("E38", "  +1.08%  "),
    @("D39", "5.427"),
    @("E39", "  -3.40%  "),
    @("D40", "0.6838"),
    @("E40", "  -2.26%  "),
    @("D41", "0.2218"),
    @("E41", "  -1.79%  "),
    @("D42", "1.252"),
    @("E42", "  -3.63%  "),
    @("D43", "8.401"),
    @("E43", "  -7.69%  "),
    @("D44", "14.53"),
    @("E44", "  -1.98%  "),
    @("D45", "0.9985"),
    @("E45", "  -0.15%  "),
    @("D46", "0.6405"),
    @("E46", "  -2.65%  "),
    @("D47", "3.876"),
    @("E47", "  -0.54%  "),
    @("D48", "2.139"),
    @("E48", "  -1.87%  "),
    @("D49", "132.71"),
    @("E49", "  -0.06%  "),
    @("D50", "0.07161"),
    @("E50", "  -1.97%  "),
    @("B51", "EOS"),
    @("C51", "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"),
    @("D51", "1.169"),
    @("E51", "  +2.30%  ")
)

foreach ($p in $pairs) {
    $addr = $p[0]
    $val = $p[1]
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}
